# Included knapsack ratio voting: refresh the per-project cost figures
# (row 2) used to compute the cost/value ratios feeding the knapsack
# voting step. Headers in row 1 (project0..project5) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2239
$ws.Range("B2").Value = 2829
$ws.Range("C2").Value = 3410
$ws.Range("D2").Value = 4074
$ws.Range("E2").Value = 5774
$ws.Range("F2").Value = 2133
